# Updated symbol list on Thu Dec 29 08:12:30 UTC 2022 with GitHub Actions
# Refreshes the coin price/volume/hour snapshot values on the active sheet,
# including a 3-row reshuffle (KickToken/BKEXToken/CEJI) in rows 41-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Numeric-looking strings must stay text (matches the sheet's inline-string cells),
    # so force a text number format while assigning, then restore the default style.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell "D2" '244.72'
Set-TextCell "G2" '8'

# Row 3
Set-TextCell "D3" '23.93'
Set-TextCell "G3" '8'

# Row 4
Set-TextCell "D4" '5.198'
Set-TextCell "G4" '8'

# Row 5
Set-TextCell "D5" '0.05724'
Set-TextCell "G5" '8'

# Row 6
Set-TextCell "D6" '6.487'
Set-TextCell "G6" '8'

# Row 7
Set-TextCell "D7" '3.170'
Set-TextCell "G7" '8'

# Row 8
Set-TextCell "D8" '0.8137'
Set-TextCell "G8" '8'

# Row 9
Set-TextCell "D9" '0.8658'
Set-TextCell "G9" '8'

# Row 10
Set-TextCell "D10" '0.1372'
Set-TextCell "G10" '8'

# Row 11
Set-TextCell "D11" '0.06942'
Set-TextCell "G11" '8'

# Row 12
Set-TextCell "D12" '0.03179'
Set-TextCell "G12" '8'

# Row 13
Set-TextCell "D13" '0.02907'
Set-TextCell "G13" '8'

# Row 14
Set-TextCell "D14" '0.09329'
Set-TextCell "G14" '8'

# Row 15
Set-TextCell "D15" '3.820'
Set-TextCell "G15" '8'

# Row 16
Set-TextCell "D16" '0.001523'
Set-TextCell "G16" '8'

# Row 17
Set-TextCell "D17" '0.04710'
Set-TextCell "G17" '8'

# Row 18
Set-TextCell "D18" '0.0005990'
Set-TextCell "G18" '8'

# Row 19
Set-TextCell "D19" '0.006150'
Set-TextCell "G19" '8'

# Row 20
Set-TextCell "D20" '0.001236'
Set-TextCell "G20" '8'

# Row 21
Set-TextCell "D21" '0.004112'
Set-TextCell "G21" '8'

# Row 22
Set-TextCell "D22" '0.00008500'
Set-TextCell "G22" '8'

# Row 23
Set-TextCell "G23" '8'

# Row 24
Set-TextCell "D24" '2.157'
Set-TextCell "G24" '8'

# Row 25
Set-TextCell "D25" '0.3195'
Set-TextCell "G25" '8'

# Row 26
Set-TextCell "G26" '8'

# Row 27
Set-TextCell "G27" '8'

# Row 28
Set-TextCell "G28" '8'

# Row 29
Set-TextCell "G29" '8'

# Row 30
Set-TextCell "G30" '8'

# Row 31
Set-TextCell "G31" '8'

# Row 32
Set-TextCell "G32" '8'

# Row 33
Set-TextCell "G33" '8'

# Row 34
Set-TextCell "G34" '8'

# Row 35
Set-TextCell "G35" '8'

# Row 36
Set-TextCell "G36" '8'

# Row 37
Set-TextCell "G37" '8'

# Row 38
Set-TextCell "G38" '8'

# Row 39
Set-TextCell "G39" '8'

# Row 40
Set-TextCell "D40" '0.03708'
Set-TextCell "G40" '8'

# Row 41
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextCell "D41" '0.006380'
$ws.Range("E41").Value = '40KickTokenKICK'
Set-TextCell "G41" '8'

# Row 42
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell "D42" '0.1054'
$ws.Range("E42").Value = '41BKEXTokenBKK'
Set-TextCell "G42" '8'

# Row 43
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextCell "D43" '0.002223'
$ws.Range("E43").Value = '42CEJICEJIWorstin24h'
Set-TextCell "G43" '8'

# Row 44
Set-TextCell "D44" '0.008109'
Set-TextCell "G44" '8'

# Row 45
Set-TextCell "D45" '0.00005447'
Set-TextCell "G45" '8'

# Row 46
Set-TextCell "G46" '8'

# Row 47
Set-TextCell "G47" '8'

# Row 48
Set-TextCell "D48" '0.002565'
Set-TextCell "G48" '8'

# Row 49
Set-TextCell "G49" '8'

# Row 50
Set-TextCell "G50" '8'

# Row 51
Set-TextCell "G51" '8'

Write-Host "Sheet refresh complete"
